$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived metrics for rows 2-13 (Resolving-Mac as sender removed; see below)
$rowUpdates = @{
    2 = @{ "G"=16.58130233333334; "H"=49.74390700000001; "I"=0.6205214785234227; "J"=0.6205214785234225; "M"=1.910418; "N"=5.731254; "O"=0.01809124304049503; "P"=0.01809124304049503; "Q"=31.677218441042; "R"=285.0949659693781; "S"=0.01122600487981455; "T"=0.01122600487981455 }
    3 = @{ "G"=16.58130233333334; "H"=49.74390700000001; "I"=0.6205214785234227; "J"=0.6205214785234225; "M"=31.995262; "N"=95.985786; "O"=0.302988173785169; "P"=0.302988173785169; "Q"=530.5231124562115; "R"=4774.708012105903; "S"=0.1880106695722848; "T"=0.1880106695722848 }
    4 = @{ "G"=16.58130233333334; "H"=49.74390700000001; "I"=0.6205214785234227; "J"=0.6205214785234225; "M"=37.858701; "N"=113.576103; "O"=0.3585136661130873; "P"=0.3585136661130873; "Q"=627.7465672282691; "R"=5649.719105054422; "S"=0.2224654301673456; "T"=0.2224654301673456 }
    5 = @{ "G"=16.58130233333334; "H"=49.74390700000001; "I"=0.6205214785234227; "J"=0.6205214785234225; "M"=33.83466466666667; "N"=101.503994; "O"=0.3204069170612486; "P"=0.3204069170612486; "Q"=561.022804184951; "R"=5049.205237664559; "S"=0.1988193739039776; "T"=0.1988193739039776 }
    6 = @{ "G"=0.5586273333333333; "H"=1.675882; "I"=0.02090549052511678; "J"=0.02090549052511678; "M"=1.910418; "N"=5.731254; "O"=0.01809124304049503; "P"=0.01809124304049503; "Q"=1.067211712892; "R"=9.604905416027998; "S"=0.0003782063099706537; "T"=0.0003782063099706537 }
    7 = @{ "G"=0.5586273333333333; "H"=1.675882; "I"=0.02090549052511678; "J"=0.02090549052511678; "M"=31.995262; "N"=95.985786; "O"=0.302988173785169; "P"=0.302988173785169; "Q"=17.87342789036133; "R"=160.860851013252; "S"=0.006334116396288287; "T"=0.006334116396288287 }
    8 = @{ "G"=0.5586273333333333; "H"=1.675882; "I"=0.02090549052511678; "J"=0.02090549052511678; "M"=37.858701; "N"=113.576103; "O"=0.3585136661130873; "P"=0.3585136661130873; "Q"=21.14890518309399; "R"=190.340146647846; "S"=0.007494904050052028; "T"=0.007494904050052028 }
    9 = @{ "G"=0.5586273333333333; "H"=1.675882; "I"=0.02090549052511678; "J"=0.02090549052511678; "M"=33.83466466666667; "N"=101.503994; "O"=0.3204069170612486; "P"=0.3204069170612486; "Q"=18.90096849696755; "R"=170.108716472708; "S"=0.006698263768805809; "T"=0.006698263768805809 }
    10 = @{ "G"=9.581631; "H"=28.744893; "I"=0.3585730309514606; "J"=0.3585730309514606; "M"=1.910418; "N"=5.731254; "O"=0.01809124304049503; "P"=0.01809124304049503; "Q"=18.304920331758; "R"=164.744282985822; "S"=0.006487031850709819; "T"=0.006487031850709819 }
    11 = @{ "G"=9.581631; "H"=28.744893; "I"=0.3585730309514606; "J"=0.3585730309514606; "M"=31.995262; "N"=95.985786; "O"=0.302988173785169; "P"=0.302988173785169; "Q"=306.566794232322; "R"=2759.101148090898; "S"=0.1086433878165959; "T"=0.1086433878165959 }
    12 = @{ "G"=9.581631; "H"=28.744893; "I"=0.3585730309514606; "J"=0.3585730309514606; "M"=37.858701; "N"=113.576103; "O"=0.3585136661130873; "P"=0.3585136661130873; "Q"=362.748103121331; "R"=3264.732928091979; "S"=0.1285533318956897; "T"=0.1285533318956897 }
    13 = @{ "G"=9.581631; "H"=28.744893; "I"=0.3585730309514606; "J"=0.3585730309514606; "M"=33.83466466666667; "N"=101.503994; "O"=0.3204069170612486; "P"=0.3204069170612486; "Q"=324.191271844738; "R"=2917.721446602642; "S"=0.1148892793884651; "T"=0.1148892793884651 }
}

foreach ($r in $rowUpdates.Keys) {
    $cols = $rowUpdates[$r]
    foreach ($colName in $cols.Keys) {
        $ws.Range("$colName$r").Value = $cols[$colName]
    }
}

# Remove the four rows for "Resolving-Mac" as the sending cluster (rows 14-17),
# which dropped out of this LR-pair export under the refreshed TPM data.
$ws.Rows("14:17").Delete()
